# Auto-generated edit script: apply scheduled market-data refresh to Sheets
# Updates cached price/profit columns (H-N) across the leve worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 5429.737
$ws.Range("I9").Value = 6420.875
$ws.Range("K9").Value = 6420.875
$ws.Range("M9").Value = -6251.875
$ws.Range("H11").Value = 44.5
$ws.Range("I11").Value = 44.5
$ws.Range("K11").Value = 44.5
$ws.Range("M11").Value = 95.5
$ws.Range("H12").Value = 13231
$ws.Range("I12").Value = 17042
$ws.Range("J12").Value = 1798
$ws.Range("K12").Value = 17042
$ws.Range("L12").Value = 1798
$ws.Range("M12").Value = -16872
$ws.Range("N12").Value = -2138
$ws.Range("H38").Value = 56.53846
$ws.Range("I38").Value = 56.53846
$ws.Range("K38").Value = 169.61538
$ws.Range("M38").Value = 202.38462
$ws.Range("H39").Value = 804.9375
$ws.Range("I39").Value = 95.916664
$ws.Range("K39").Value = 287.749992
$ws.Range("M39").Value = 8.25000799999998
$ws.Range("H69").Value = 14899.2
$ws.Range("I69").Value = 9872.75
$ws.Range("J69").Value = 16727
$ws.Range("K69").Value = 29618.25
$ws.Range("L69").Value = 50181
$ws.Range("M69").Value = -28744.25
$ws.Range("N69").Value = -51929
$ws.Range("H70").Value = 14272.637
$ws.Range("J70").Value = 25750
$ws.Range("L70").Value = 77250
$ws.Range("N70").Value = -77790
$ws.Range("H72").Value = 14899.2
$ws.Range("I72").Value = 9872.75
$ws.Range("J72").Value = 16727
$ws.Range("K72").Value = 88854.75
$ws.Range("L72").Value = 150543
$ws.Range("M72").Value = -84486.75
$ws.Range("N72").Value = -159279
$ws.Range("H73").Value = 14272.637
$ws.Range("J73").Value = 25750
$ws.Range("L73").Value = 77250
$ws.Range("N73").Value = -79122
$ws.Range("H113").Value = 3082.6365
$ws.Range("J113").Value = 4802
$ws.Range("L113").Value = 4802
$ws.Range("N113").Value = -11310
$ws.Range("H137").Value = 3489.7083
$ws.Range("I137").Value = 3083
$ws.Range("J137").Value = 4303.125
$ws.Range("K137").Value = 9249
$ws.Range("L137").Value = 12909.375
$ws.Range("M137").Value = -6699
$ws.Range("N137").Value = -18009.375
$ws.Range("H138").Value = 4915.8125
$ws.Range("I138").Value = 4335.8887
$ws.Range("J138").Value = 4975.8047
$ws.Range("K138").Value = 13007.6661
$ws.Range("L138").Value = 14927.4141
$ws.Range("M138").Value = -7867.666100000002
$ws.Range("N138").Value = -25207.4141

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20062.396
$ws.Range("I32").Value = 14802.552
$ws.Range("J32").Value = 28090.578
$ws.Range("K32").Value = 14802.552
$ws.Range("L32").Value = 28090.578
$ws.Range("M32").Value = -14515.552
$ws.Range("N32").Value = -28664.578
$ws.Range("H61").Value = 58827910
$ws.Range("J61").Value = 4849.8
$ws.Range("L61").Value = 4849.8
$ws.Range("N61").Value = -5273.8
$ws.Range("H102").Value = 1853228.9
$ws.Range("I102").Value = 1950635.9
$ws.Range("K102").Value = 1950635.9
$ws.Range("M102").Value = -1949013.9
$ws.Range("H132").Value = 3040855.2
$ws.Range("I132").Value = 3135808
$ws.Range("K132").Value = 9407424
$ws.Range("M132").Value = -9404894
$ws.Range("H136").Value = 58827910
$ws.Range("J136").Value = 4849.8
$ws.Range("L136").Value = 14549.4
$ws.Range("N136").Value = -19649.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1424.75
$ws.Range("I99").Value = 1359.7
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 1359.7
$ws.Range("L99").Value = 1750
$ws.Range("M99").Value = 138.3
$ws.Range("N99").Value = -4746
$ws.Range("H107").Value = 31338.334
$ws.Range("I107").Value = 1154.7858
$ws.Range("K107").Value = 1154.7858
$ws.Range("M107").Value = 765.2141999999999
$ws.Range("H134").Value = 17245524
$ws.Range("I134").Value = 20003450
$ws.Range("J134").Value = 8497.25
$ws.Range("K134").Value = 60010350
$ws.Range("L134").Value = 25491.75
$ws.Range("M134").Value = -60007815
$ws.Range("N134").Value = -30561.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7387.4116
$ws.Range("I31").Value = 4489.645
$ws.Range("K31").Value = 4489.645
$ws.Range("M31").Value = -4194.645
$ws.Range("H34").Value = 7387.4116
$ws.Range("I34").Value = 4489.645
$ws.Range("K34").Value = 4489.645
$ws.Range("M34").Value = -4287.645
$ws.Range("H58").Value = 22733858
$ws.Range("I58").Value = 29419620
$ws.Range("J58").Value = 2266.4
$ws.Range("K58").Value = 29419620
$ws.Range("L58").Value = 2266.4
$ws.Range("M58").Value = -29419417
$ws.Range("N58").Value = -2672.4
$ws.Range("H105").Value = 3573383.5
$ws.Range("I105").Value = 5953763.5
$ws.Range("K105").Value = 5953763.5
$ws.Range("M105").Value = -5952016.5
$ws.Range("H136").Value = 22733858
$ws.Range("I136").Value = 29419620
$ws.Range("J136").Value = 2266.4
$ws.Range("K136").Value = 88258860
$ws.Range("L136").Value = 6799.200000000001
$ws.Range("M136").Value = -88256310
$ws.Range("N136").Value = -11899.2
$ws.Range("H141").Value = 221990.12
$ws.Range("J141").Value = 293709.66
$ws.Range("L141").Value = 293709.66
$ws.Range("N141").Value = -304069.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1986.3
$ws.Range("I122").Value = 1284.4445
$ws.Range("J122").Value = 2560.5454
$ws.Range("K122").Value = 11560.0005
$ws.Range("L122").Value = 23044.9086
$ws.Range("M122").Value = -9110.0005
$ws.Range("N122").Value = -27944.9086
$ws.Range("J137").Value = 3498
$ws.Range("L137").Value = 10494
$ws.Range("N137").Value = -20694
$ws.Range("H141").Value = 1669112.5
$ws.Range("I141").Value = 1669112.5
$ws.Range("K141").Value = 5007337.5
$ws.Range("M141").Value = -5002157.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 60205
$ws.Range("I113").Value = 114803.336
$ws.Range("K113").Value = 114803.336
$ws.Range("M113").Value = -112633.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1820.9333
$ws.Range("I22").Value = 1786.5
$ws.Range("K22").Value = 1786.5
$ws.Range("M22").Value = -1491.5
$ws.Range("H27").Value = 1820.9333
$ws.Range("I27").Value = 1786.5
$ws.Range("K27").Value = 1786.5
$ws.Range("M27").Value = -1679.5
$ws.Range("H55").Value = 612.0345
$ws.Range("I55").Value = 358.13333
$ws.Range("J55").Value = 884.0714
$ws.Range("K55").Value = 358.13333
$ws.Range("L55").Value = 884.0714
$ws.Range("M55").Value = -185.13333
$ws.Range("N55").Value = -1230.0714
$ws.Range("H93").Value = 1762.2222
$ws.Range("J93").Value = 2001.6
$ws.Range("L93").Value = 2001.6
$ws.Range("N93").Value = -4497.6
$ws.Range("H132").Value = 20888336
$ws.Range("I132").Value = 22193562
$ws.Range("K132").Value = 66580686
$ws.Range("M132").Value = -66578156
$ws.Range("H136").Value = 2704.5293
$ws.Range("I136").Value = 2686.125
$ws.Range("K136").Value = 8058.375
$ws.Range("M136").Value = -5508.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13853.714
$ws.Range("J45").Value = 13853.714
$ws.Range("L45").Value = 13853.714
$ws.Range("N45").Value = -14835.714
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("H107").Value = 938.619
$ws.Range("I107").Value = 435.44446
$ws.Range("J107").Value = 1316
$ws.Range("K107").Value = 1306.33338
$ws.Range("L107").Value = 3948
$ws.Range("M107").Value = 613.66662
$ws.Range("N107").Value = -7788

# Cells whose profit value became unavailable (removed entirely) after refresh
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

